$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2215189873417721
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.02848101265822785
$ws.Range("P2").Value = 0.1424050632911392
$ws.Range("S2").Value = 0.1075949367088608
$ws.Range("B3").Value = 0.006097560975609756
$ws.Range("C3").Value = 0.03658536585365853
$ws.Range("J3").Value = 0.01829268292682927
$ws.Range("P3").Value = 0.7073170731707317
$ws.Range("S3").Value = 0.2317073170731707
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.6444444444444445
$ws.Range("S4").Value = 0.3111111111111111
$ws.Range("B6").Value = 0.05128205128205128
$ws.Range("D6").Value = 0.004273504273504274
$ws.Range("E6").Value = 0.004273504273504274
$ws.Range("F6").Value = 0.07264957264957266
$ws.Range("J6").Value = 0.2606837606837607
$ws.Range("O6").Value = 0.01282051282051282
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.07692307692307693
$ws.Range("S6").Value = 0.3632478632478632
$ws.Range("B7").Value = 0.08205128205128205
$ws.Range("D7").Value = 0.02564102564102564
$ws.Range("E7").Value = 0.005128205128205128
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.1897435897435897
$ws.Range("O7").Value = 0.02051282051282051
$ws.Range("Q7").Value = 0.1128205128205128
$ws.Range("R7").Value = 0.1230769230769231
$ws.Range("S7").Value = 0.3743589743589744
$ws.Range("B8").Value = 0.07644628099173553
$ws.Range("D8").Value = 0.004132231404958678
$ws.Range("F8").Value = 0.05578512396694215
$ws.Range("J8").Value = 0.09504132231404959
$ws.Range("O8").Value = 0.006198347107438017
$ws.Range("Q8").Value = 0.1900826446280992
$ws.Range("R8").Value = 0.1198347107438017
$ws.Range("S8").Value = 0.4524793388429752
$ws.Range("B9").Value = 0.1182266009852217
$ws.Range("D9").Value = 0.01477832512315271
$ws.Range("F9").Value = 0.06403940886699508
$ws.Range("J9").Value = 0.1182266009852217
$ws.Range("O9").Value = 0.004926108374384237
$ws.Range("Q9").Value = 0.1527093596059113
$ws.Range("R9").Value = 0.06896551724137931
$ws.Range("S9").Value = 0.458128078817734
$ws.Range("B10").Value = 0.1094674556213018
$ws.Range("D10").Value = 0.02662721893491124
$ws.Range("F10").Value = 0.07100591715976332
$ws.Range("J10").Value = 0.136094674556213
$ws.Range("O10").Value = 0.009615384615384616
$ws.Range("Q10").Value = 0.1871301775147929
$ws.Range("R10").Value = 0.08579881656804733
$ws.Range("S10").Value = 0.3742603550295858
$ws.Range("G11").Value = 0.1757188498402556
$ws.Range("J11").Value = 0.07348242811501597
$ws.Range("K11").Value = 0.2332268370607029
$ws.Range("L11").Value = 0.5079872204472844
$ws.Range("S11").Value = 0.009584664536741214
$ws.Range("G12").Value = 0.7407407407407407
$ws.Range("J12").Value = 0.228395061728395
$ws.Range("L12").Value = 0.006172839506172839
$ws.Range("S12").Value = 0.02469135802469136
$ws.Range("G13").Value = 0.5581395348837209
$ws.Range("J13").Value = 0.3953488372093023
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.01052631578947368
$ws.Range("H15").Value = 0.1894736842105263
$ws.Range("I15").Value = 0.08421052631578947
$ws.Range("J15").Value = 0.3789473684210526
$ws.Range("K15").Value = 0.05263157894736842
$ws.Range("M15").Value = 0.01052631578947368
$ws.Range("O15").Value = 0.1
$ws.Range("S15").Value = 0.1736842105263158
$ws.Range("F16").Value = 0.00546448087431694
$ws.Range("H16").Value = 0.2131147540983606
$ws.Range("I16").Value = 0.08743169398907104
$ws.Range("J16").Value = 0.4098360655737705
$ws.Range("K16").Value = 0.09836065573770492
$ws.Range("M16").Value = 0.01092896174863388
$ws.Range("O16").Value = 0.03278688524590164
$ws.Range("S16").Value = 0.1420765027322404
$ws.Range("F17").Value = 0.01624129930394431
$ws.Range("H17").Value = 0.2018561484918794
$ws.Range("I17").Value = 0.1020881670533643
$ws.Range("J17").Value = 0.4106728538283063
$ws.Range("K17").Value = 0.09280742459396751
$ws.Range("M17").Value = 0.0185614849187935
$ws.Range("N17").Value = 0.002320185614849188
$ws.Range("O17").Value = 0.04872389791183294
$ws.Range("S17").Value = 0.1067285382830626
$ws.Range("F18").Value = 0.01739130434782609
$ws.Range("H18").Value = 0.1826086956521739
$ws.Range("I18").Value = 0.08260869565217391
$ws.Range("J18").Value = 0.4217391304347826
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.02608695652173913
$ws.Range("O18").Value = 0.06521739130434782
$ws.Range("S18").Value = 0.1
$ws.Range("F19").Value = 0.018140589569161
$ws.Range("H19").Value = 0.2169312169312169
$ws.Range("I19").Value = 0.08163265306122448
$ws.Range("J19").Value = 0.3801965230536659
$ws.Range("K19").Value = 0.1111111111111111
$ws.Range("M19").Value = 0.01965230536659108
$ws.Range("N19").Value = 0.001511715797430083
$ws.Range("O19").Value = 0.06198034769463341
$ws.Range("S19").Value = 0.108843537414966
